$d = $word.ActiveDocument

# 1) Insert a new bullet paragraph "(Familia, hermano, padre, tio);" right
#    before the "(etc.: CSPO, Kinds, Statements)..." bullet. The new
#    paragraph inherits the list/paragraph formatting of the paragraph it
#    is inserted before.
$r = $d.Content
$found = $r.Find.Execute("(etc.: CSPO, Kinds, Statements). ResourceOccurrences LHS, Concepts (ResourceOccurrence Context Kind), RHS:", $true)
if ($found) {
    $targetPara = $r.Paragraphs(1)
    $idx = $targetPara.Index
    $targetPara.Range.InsertParagraphBefore()
    $d.Paragraphs($idx).Range.Text = "(Familia, hermano, padre, tio);"
}

# 2) Replace "(Kind, Kind, CSPO)" with "(Kind, Kind, Kind)".
$d.Content.Find.Execute("(Kind, Kind, CSPO)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "(Kind, Kind, Kind)", 2)

# 3) Insert a new bullet paragraph with the "Rules Aggregation: ..." text
#    right after the "(CSPO, CSPO, Kind)" bullet.
$r2 = $d.Content
$found2 = $r2.Find.Execute("(CSPO, CSPO, Kind)", $true)
if ($found2) {
    $srcPara = $r2.Paragraphs(1)
    $nextPara = $srcPara.Next()
    $nidx = $nextPara.Index
    $nextPara.Range.InsertParagraphBefore()
    $d.Paragraphs($nidx).Range.Text = "Rules Aggregation: Rule Context application matches / filters input Statements Flux for LHS Statements, Kinds, CSPOs Flux. Concepts Aggregated by Context CSPOs, Kinds, Statements. RHS result of applying Concept Kind to LHS Flux (infer Grammar)."
}
